$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("F405_pill-top-pos")

# C3's value (row 4, column B) changes from "18p" to "12p".
# Because the same text "18p" is also used by C6 (row 7, column B) via the
# shared string table, updating the shared string text updates both cells.
$ws.Range("B4").Value = "12p"
$ws.Range("B7").Value = "12p"

# R9's value (row 32, column B) changes from "10k" to "20k".
$ws.Range("B32").Value = "20k"

# U1's package (row 35, column C) changes from "my_STM32F405RGTx_2" to
# "LQFP-64_10x10mm_P0.5mm".
$ws.Range("C35").Value = "LQFP-64_10x10mm_P0.5mm"

# Reflect the final selected cell recorded in the saved workbook.
$ws.Range("I12").Select()
